$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = "السجل"
$ws.Range("B7").Value = "السجل"

$ws.Range("A8").Value = "الاسم"
$ws.Range("B8").Value = "الاسم"

$ws.Range("A9").Value = "العنوان"
$ws.Range("B9").Value = "العنوان"

$ws.Range("B9").Select()
